$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the username e-mail shown in A2 (sharedString text edit) ---
# "sanitycheck9250@yopmail.com" -> "sanitycheck9150@yopmail.com"
$ws.Range("A2").Value = "sanitycheck9150@yopmail.com"

# --- 2. Swap which cell each hyperlink points at -----------------------
# Originally: B2 -> mailto:Tesh@1234 , A2 -> mailto:sanitycheck9250@yopmail.com
# After:      A2 -> mailto:Tesh@1234 , B2 -> mailto:sanitycheck9250@yopmail.com
# (the link target text itself is untouched by this edit, only the cell it
#  is attached to changes)
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:Tesh@1234")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:sanitycheck9250@yopmail.com")

# Re-adding the hyperlinks resets the cell formatting on A2/B2; put the
# original "Hyperlink" cell style back so the underlying cell style index
# used by the sheet is unchanged.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("B2").Style = "Hyperlink"
